# Generate Report for Handoff
# Updates the localization-status workbook: file #1 (2f9b3e4b...) is replaced by
# d8671e60-fb43-46ec-b607-b1158c06422d (status "Ready for handoff"), and file #2
# (857049c3...) is replaced by ffff59b0e27c-8ce4-4f8e-9c03-8b6f4dc11a6e, which now
# shares the same handoff target/file info as file #1.

$wb = $excel.ActiveWorkbook

$guid1 = "d8671e60-fb43-46ec-b607-b1158c06422d"
$guid2 = "ffff59b0e27c-8ce4-4f8e-9c03-8b6f4dc11a6e"
$hash1 = "81bf13d9d024ef5e179a89dfe1e41c64b4cf54d0"

$status = "Ready for handoff"
$overviewDate = "2016-16-13 23:16:07"
$zhHandoffDate = "2016-03-13 23:16:03"
$deHandoffDate = "2016-03-13 23:16:07"
$handbackDate = "0001-01-01 00:00:00"

$zhXlf = "$guid1.$hash1.zh-cn.xlf"
$deXlf = "$guid1.$hash1.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "$guid1.md"
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $overviewDate

$ws1.Range("A3").Value = "$guid2.md"
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $overviewDate

$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid1.md", "", "", "$guid1.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid2.md", "", "", "$guid2.md")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "$guid1.md"
$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $zhXlf
$ws2.Range("E2").Value = $zhHandoffDate
$ws2.Range("H2").Value = $handbackDate

$ws2.Range("A3").Value = "$guid2.md"
$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $zhXlf
$ws2.Range("E3").Value = $zhHandoffDate
$ws2.Range("H3").Value = $handbackDate

$ws2.Range("F2:G3").Clear()

$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid1.md", "", "", "$guid1.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid1.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9b777bfc22420acedd49d95fa1cd066342afde/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid2.md", "", "", "$guid2.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid2.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9b777bfc22420acedd49d95fa1cd066342afde/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "$guid1.md"
$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $deXlf
$ws3.Range("E2").Value = $deHandoffDate
$ws3.Range("H2").Value = $handbackDate

$ws3.Range("A3").Value = "$guid2.md"
$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $deXlf
$ws3.Range("E3").Value = $deHandoffDate
$ws3.Range("H3").Value = $handbackDate

$ws3.Range("F2:G3").Clear()

$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid1.md", "", "", "$guid1.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid1.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/530973197b39cc3e4148f5f5327995001a6ff49b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid2.md", "", "", "$guid2.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e2f4d693f76102239c2e43c768cdba4b770f1a9/e2e/$guid2.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/530973197b39cc3e4148f5f5327995001a6ff49b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf)

Write-Host "Report for Handoff generated."
